$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 45133
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 80
$ws.Range("K2").Value = 2500
$ws.Range("L2").Value = 2500
$ws.Range("M2").Value = 2500
$ws.Range("P2").Value = 833

$ws.Range("D3").Value = 45191
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 100
$ws.Range("K3").Value = 2500
$ws.Range("L3").Value = 2500
$ws.Range("M3").Value = 2500
$ws.Range("P3").Value = 833

$ws.Range("D4").Value = 45166
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 2500
$ws.Range("L4").Value = 2500
$ws.Range("M4").Value = 2500
$ws.Range("P4").Value = 833

$ws.Range("D5").Value = 45160
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 2500
$ws.Range("L5").Value = 2500
$ws.Range("M5").Value = 2500
$ws.Range("P5").Value = 833

$ws.Range("D6").Value = 45145
$ws.Range("I6").Value = "Primera"
$ws.Range("J6").Value = 60
$ws.Range("K6").Value = 2500
$ws.Range("L6").Value = 2500
$ws.Range("M6").Value = 2500
$ws.Range("P6").Value = 833

$ws.Range("D7").Value = 45145
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 80
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 2000
$ws.Range("M7").Value = 2000
$ws.Range("P7").Value = 667

$ws.Range("D8").Value = 45148
$ws.Range("I8").Value = "Primera"
$ws.Range("J8").Value = 80
$ws.Range("K8").Value = 2500
$ws.Range("L8").Value = 2500
$ws.Range("M8").Value = 2500
$ws.Range("P8").Value = 833

$ws.Range("D9").Value = 45148
$ws.Range("I9").Value = "Segunda"
$ws.Range("J9").Value = 60
$ws.Range("K9").Value = 2000
$ws.Range("L9").Value = 2000
$ws.Range("M9").Value = 2000
$ws.Range("P9").Value = 667

$ws.Range("D10").Value = 45173
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 100
$ws.Range("K10").Value = 2500
$ws.Range("L10").Value = 2500
$ws.Range("M10").Value = 2500
$ws.Range("P10").Value = 833

$ws.Range("D11").Value = 44846
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 200
$ws.Range("K11").Value = 1200
$ws.Range("L11").Value = 1300
$ws.Range("M11").Value = 1250
$ws.Range("P11").Value = 417

$ws.Range("D12").Value = 44846
$ws.Range("I12").Value = "Segunda"
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 1000
$ws.Range("L12").Value = 1000
$ws.Range("M12").Value = 1000
$ws.Range("P12").Value = 333

$ws.Range("D13").Value = 45161
$ws.Range("I13").Value = "Primera"
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 2500
$ws.Range("L13").Value = 2500
$ws.Range("M13").Value = 2500
$ws.Range("P13").Value = 833

$ws.Range("D14").Value = 45146
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 80
$ws.Range("K14").Value = 2500
$ws.Range("L14").Value = 2500
$ws.Range("M14").Value = 2500
$ws.Range("P14").Value = 833

$ws.Range("D15").Value = 45146
$ws.Range("I15").Value = "Segunda"
$ws.Range("J15").Value = 80
$ws.Range("K15").Value = 2000
$ws.Range("L15").Value = 2000
$ws.Range("M15").Value = 2000
$ws.Range("P15").Value = 667

$ws.Range("D16").Value = 45163
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 100
$ws.Range("K16").Value = 2500
$ws.Range("L16").Value = 2500
$ws.Range("M16").Value = 2500
$ws.Range("P16").Value = 833

$ws.Range("D17").Value = 45195
$ws.Range("I17").Value = "Primera"
$ws.Range("J17").Value = 100
$ws.Range("K17").Value = 2500
$ws.Range("L17").Value = 2500
$ws.Range("M17").Value = 2500
$ws.Range("P17").Value = 833

$ws.Range("D18").Value = 44838
$ws.Range("I18").Value = "Primera"
$ws.Range("J18").Value = 200
$ws.Range("K18").Value = 1200
$ws.Range("L18").Value = 1300
$ws.Range("M18").Value = 1250
$ws.Range("P18").Value = 417

$ws.Range("D19").Value = 44838
$ws.Range("I19").Value = "Segunda"
$ws.Range("J19").Value = 150
$ws.Range("K19").Value = 1000
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = 1000
$ws.Range("P19").Value = 333

$ws.Range("D20").Value = 45135
$ws.Range("I20").Value = "Primera"
$ws.Range("J20").Value = 70
$ws.Range("K20").Value = 2500
$ws.Range("L20").Value = 2500
$ws.Range("M20").Value = 2500
$ws.Range("P20").Value = 833

$ws.Range("D21").Value = 45149
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 80
$ws.Range("K21").Value = 2500
$ws.Range("L21").Value = 2500
$ws.Range("M21").Value = 2500
$ws.Range("P21").Value = 833

$ws.Range("D22").Value = 45149
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 80
$ws.Range("K22").Value = 2000
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = 2000
$ws.Range("P22").Value = 667

$ws.Range("D23").Value = 45134
$ws.Range("I23").Value = "Primera"
$ws.Range("J23").Value = 50
$ws.Range("K23").Value = 2500
$ws.Range("L23").Value = 2500
$ws.Range("M23").Value = 2500
$ws.Range("P23").Value = 833

$ws.Range("D24").Value = 45175
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 150
$ws.Range("K24").Value = 2500
$ws.Range("L24").Value = 2500
$ws.Range("M24").Value = 2500
$ws.Range("P24").Value = 833

$ws.Range("D25").Value = 45176
$ws.Range("I25").Value = "Primera"
$ws.Range("J25").Value = 100
$ws.Range("K25").Value = 2500
$ws.Range("L25").Value = 2500
$ws.Range("M25").Value = 2500
$ws.Range("P25").Value = 833

$ws.Range("D26").Value = 44832
$ws.Range("I26").Value = "Primera"
$ws.Range("J26").Value = 200
$ws.Range("K26").Value = 1200
$ws.Range("L26").Value = 1300
$ws.Range("M26").Value = 1250
$ws.Range("P26").Value = 417

$ws.Range("D27").Value = 44832
$ws.Range("I27").Value = "Segunda"
$ws.Range("J27").Value = 150
$ws.Range("K27").Value = 1000
$ws.Range("L27").Value = 1000
$ws.Range("M27").Value = 1000
$ws.Range("P27").Value = 333
